# Update cryptocurrency price/volume figures in the worksheet.
# Source: cryptos list refresh commit (GitHub Actions automated update).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (Price column D, Volume(1h) column E).
$updates = @{
    "D2" = '26.461.86'
    "E2" = '  +1.53%  '
    "D3" = '1.676.83'
    "E3" = '  +2.54%  '
    "E4" = '  -0.09%  '
    "D5" = '217.39'
    "E5" = '  +2.03%  '
    "D6" = '0.5315'
    "E6" = '  +1.56%  '
    "E7" = '  -0.09%  '
    "D8" = '0.2693'
    "E8" = '  +3.98%  '
    "D9" = '0.06415'
    "E9" = '  +2.07%  '
    "D10" = '21.80'
    "E10" = '  +5.42%  '
    "D11" = '0.07829'
    "E11" = '  +2.26%  '
    "D12" = '1.678.13'
    "E12" = '  +2.45%  '
    "D13" = '4.507'
    "E13" = '  +2.21%  '
    "D14" = '0.5584'
    "E14" = '  +1.53%  '
    "D15" = '0.0₅8346'
    "E15" = '  +2.25%  '
    "D16" = '65.77'
    "E16" = '  +1.43%  '
    "D17" = '26.507.16'
    "E17" = '  +1.74%  '
    "E18" = '  -0.07%  '
    "D19" = '4.743'
    "E19" = '  +1.30%  '
    "D20" = '194.34'
    "E20" = '  +3.29%  '
    "D21" = '10.30'
    "E21" = '  +1.51%  '
    "E22" = '  +3.56%  '
    "E23" = '  -0.07%  '
    "D24" = '142.41'
    "E24" = '  -2.20%  '
    "D25" = '0.1289'
    "D26" = '7.396'
    "E26" = '  +0.07%  '
    "E27" = '  +3.31%  '
    "E28" = '  +3.21%  '
    "D29" = '0.06310'
    "E29" = '  +6.09%  '
    "D30" = '1.277'
    "E30" = '  +1.82%  '
    "D31" = '3.643'
    "E31" = '  +6.01%  '
    "E32" = '  +1.58%  '
    "D33" = '1.679'
    "E33" = '  +2.83%  '
    "E34" = '  +2.93%  '
    "D35" = '0.6205'
    "E35" = '  +8.74%  '
    "D36" = '2.420'
    "E36" = '  +1.01%  '
    "D37" = '2.789'
    "E37" = '  +0.98%  '
    "D38" = '0.01638'
    "E38" = '  +1.46%  '
    "D39" = '6.162'
    "E39" = '  +7.47%  '
    "D40" = '1.092.51'
    "E40" = '  +5.81%  '
    "D41" = '0.8659'
    "E41" = '  +1.72%  '
    "E43" = '  +0.18%  '
    "D44" = '1.822.73'
    "E44" = '  +2.08%  '
    "D45" = '57.74'
    "E45" = '  +3.44%  '
    "D46" = '8.190'
    "E46" = '  +2.18%  '
    "E47" = '  +0.01%  '
    "E48" = '  -1.27%  '
    "D49" = '1.492'
    "E49" = '  +7.91%  '
    "D50" = '0.05208'
    "E50" = '  +0.84%  '
    "E51" = '  +2.59%  '
}

# Regex used to detect values that Excel would otherwise auto-convert into a
# numeric type (e.g. "217.39"). Such cells must be written as literal text,
# matching the original inlineStr content, so we temporarily force a text
# number format while assigning the value and then restore the default
# "Normal" style straight after (leaving the cell's formatting untouched).
$numericPattern = '^-?\d+(\.\d+)?$'

foreach ($address in $updates.Keys) {
    $value = $updates[$address]
    $range = $ws.Range($address)

    if ($value -match $numericPattern) {
        $range.NumberFormat = "@"
        $range.Value = $value
        $range.Style = "Normal"
    } else {
        $range.Value = $value
    }
}
